$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B/C/D/E hold plain text in this sheet (t="inlineStr" originally).
# For Price values (column D) that look like plain decimal numbers, Excel
# would silently convert the assigned text into a real number. Prefix
# those with an apostrophe to force text, then reset the cell style back
# to Normal so the "quote prefix" formatting flag is not left behind.

$ws.Range("D2").Value = "41.939.85"
$ws.Range("E2").Value = "  -1.15%  "

$ws.Range("D3").Value = "2.215.94"
$ws.Range("E3").Value = "  -1.59%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'241.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.24%  "

$ws.Range("D6").Value = "'0.619"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.40%  "

$ws.Range("D7").Value = "'72.94"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.89%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.604"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.88%  "

$ws.Range("D10").Value = "'42.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.75%  "

$ws.Range("E11").Value = "  +0.61%  "

$ws.Range("D12").Value = "'7.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.84%  "

$ws.Range("E13").Value = "  +0.34%  "

$ws.Range("D14").Value = "2.552.42"
$ws.Range("E14").Value = "  -1.36%  "

$ws.Range("E15").Value = "  -2.69%  "

$ws.Range("D16").Value = "'0.835"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.70%  "

$ws.Range("D17").Value = "2.215.40"
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").Value = "41.851.92"
$ws.Range("E18").Value = "  -0.83%  "

$ws.Range("E19").Value = "  +4.93%  "

$ws.Range("D20").Value = "'6.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("D21").Value = "'72.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").Value = "'11.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +20.65%  "

$ws.Range("D23").Value = "'229.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.94%  "

$ws.Range("E24").Value = "  -7.22%  "

$ws.Range("D25").Value = "'11.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.45%  "

$ws.Range("D27").Value = "'3.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.25%  "

$ws.Range("E28").Value = "  -1.54%  "

$ws.Range("E29").Value = "  -1.40%  "

$ws.Range("D30").Value = "'167.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.81%  "

$ws.Range("D31").Value = "'20.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.19%  "

$ws.Range("D32").Value = "'5.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.37%  "

$ws.Range("D33").Value = "'0.0799"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.55%  "

$ws.Range("D34").Value = "'30.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.56%  "

$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").Value = "'0.109"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.97%  "

$ws.Range("D37").Value = "'4.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.91%  "

$ws.Range("D38").Value = "'0.0302"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.45%  "

$ws.Range("D39").Value = "'13.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.05%  "

$ws.Range("D40").Value = "'65.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.16%  "

$ws.Range("E41").Value = "  -3.27%  "

$ws.Range("E42").Value = "  -2.98%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.198"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'8.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.35%  "

$ws.Range("D45").Value = "'104.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.11%  "

$ws.Range("D46").Value = "'0.100"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.39%  "

$ws.Range("D47").Value = "'2.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.31%  "

$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "'1.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.38%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'1.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.45%  "

$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("D51").Value = "2.424.31"
$ws.Range("E51").Value = "  -1.57%  "
